$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking Price strings stay as text (matching the
# original inline-string cell type) instead of Excel auto-converting
# them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.345.77"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.94"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.19"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6289"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07430"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.84"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.843.19"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.982"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6781"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001022"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.02"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.248"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.349.78"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.26"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.437"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9990"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.471"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1355"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.42"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06512"
$ws.Range("E28").Value = "  +14.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.450"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.487"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.065"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.063"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6951"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01856"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.815"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.239.43"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.782"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9321"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.990.44"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.84"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.64"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.053"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.714"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1153"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.973"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3903"
$ws.Range("E51").Value = "  -1.70%  "
